$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2-3: account holder name / card number ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long, purely-numeric account number that must stay TEXT (it was
# an inline string before the edit). Assigning a digit-only string directly
# makes Excel coerce it to a number, so force a temporary Text number format,
# write the value, then restore the original (General) formatting by pasting
# the formats from a neighbouring cell that already carries the target style
# (C3 uses the same style as B3) - this avoids inventing a brand-new style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Copy() | Out-Null
$ws.Range("B3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("C3").Value = "Mohaupt"

# --- Row 5: opening balance date line ---
$ws.Range("D5").Value = "KONTOSTAND AM 04.03.2024"

# --- Row 6: first transaction ---
$ws.Range("B6").Value = "07.03."
$ws.Range("C6").Value = "08.03."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 83073115"
$ws.Range("E6").Value = "85,68-"

# --- Row 7: second transaction ---
$ws.Range("B7").Value = "11.03."
$ws.Range("C7").Value = "12.03."
$ws.Range("D7").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E7").Value = "24,64-"

# --- Row 8: third transaction ---
$ws.Range("B8").Value = "12.03."
$ws.Range("C8").Value = "13.03."
$ws.Range("D8").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E8").Value = "81,22-"

# --- Row 9: new fourth transaction (previously a blank filler row) ---
$ws.Range("B9").Value = "13.03."
$ws.Range("C9").Value = "14.03."
$ws.Range("D9").Value = "AMAZON.DE MKTPLC EU ESZWSW"
$ws.Range("E9").Value = "249,66-"
# E9 switches from the centered "blank filler" style to the right-aligned
# amount style already used by E6:E8/E12. Copy that exact formatting across
# instead of poking individual alignment properties so no new style is added.
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Row 12: closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 17.03.2024"
$ws.Range("E12").Value = "441,20-"

# --- Row 13: next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 26.03.2024"

$excel.CutCopyMode = $false
